$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed login/email values (A4, A5, A6).
# The other cells in the data table are unaffected by this commit.
$ws.Range("A5").Value = "ahmet@yahoo.com"
$ws.Range("A6").Value = "firat@gmail.com"
$ws.Range("A4").Value = "burcu@yahoo.com"

# Move the active selection from C6 to A4, matching the saved view state.
$ws.Range("A4").Select()
